# Apply the TeSt1 Anthropometrics update:
#  - correct several existing dates / a Day_Type value
#  - append a new row (7) of data with a new Weight (H) reading
#  - leave the selection on the newly entered cell, matching the
#    author's last cursor position when they saved the workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix previously mis-entered dates / Day_Type -------------------------
$ws.Range("B2").Value = 43914
$ws.Range("B3").Value = 43946
$ws.Range("B4").Value = 43977
$ws.Range("C4").Value = 3
$ws.Range("B6").Value = 44040

# --- append the new measurement row --------------------------------------
# Copy an existing date cell first so the new date cell (B7) inherits the
# same short-date number format style already used by the column, then
# overwrite it with the real value.
$ws.Range("B6").Copy($ws.Range("B7"))

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 44045
$ws.Range("C7").Value = 1
$ws.Range("H7").Value = 62.3

# Column B widened a touch to keep fitting the date text after the edits
# (mirrors Excel's "best fit" recalculation once the new date was entered).
$ws.Columns.Item(2).ColumnWidth = 8.83

# Reflect the author's final selection/cursor position.
$ws.Range("H7").Select()
